$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The SpritePNGtoCubePixelHelper.exe + textures now sit one folder deeper
# relative to this workbook (it moved into PLY_Models\Defender\), so both
# the exe path and the PNG (relative) location need one more "../".
# Leading "'" mirrors how these were originally typed in Excel (forces
# text / keeps the existing quote-prefix cell formatting).
$ws.Range("B1").Value = "'../../../x64/Release/SpritePNGtoCubePixelHelper.exe"
$ws.Range("B2").Value = "'../../textures/Defender/Sprites (isolated)"

# Reset the view: scroll back to the top, zoom in, and leave a single
# cell selected instead of the big B6:B91 range.
$ws.Range("B3").Select()
$excel.ActiveWindow.Zoom = 130
